$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$data = @(
    @{Row=2;  B=2450.05;  C=2419.25},
    @{Row=3;  B=384.75;   C=388.4},
    @{Row=4;  B=1514.5;   C=1488.9},
    @{Row=5;  B=7320.6;   C=7163},
    @{Row=6;  B=235.1;    C=229.6},
    @{Row=7;  B=190.6;    C=187.15},
    @{Row=8;  B=44601.5;  C=43989.15},
    @{Row=9;  B=503.25;   C=504.1},
    @{Row=10; B=3400.4;   C=3336.6},
    @{Row=11; B=143.8;    C=143.95},
    @{Row=12; B=1185.8;   C=1172.15},
    @{Row=13; B=1416.55;  C=1378.15},
    @{Row=14; B=701.1;    C=684.5},
    @{Row=15; B=433.2;    C=423.35},
    @{Row=16; B=1591.95;  C=1575.4},
    @{Row=17; B=299.15;   C=297.65},
    @{Row=18; B=19512.1;  C=19253.8},
    @{Row=19; B=569.55;   C=561.35},
    @{Row=20; B=620.8;    C=617},
    @{Row=21; B=611.2;    C=601},
    @{Row=22; B=255.35;   C=245.1},
    @{Row=23; B=127.05;   C=122.9}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
    $ws.Cells.Item($entry.Row, 3).Value = $entry.C
}
